$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.839.44"
$ws.Range("E2").Value = "'  -1.86%  "
$ws.Range("D3").Value = "'1.631.84"
$ws.Range("E3").Value = "'  -2.06%  "
$ws.Range("D4").Value = "'1.009"
$ws.Range("E4").Value = "'  +0.05%  "
$ws.Range("D5").Value = "'215.17"
$ws.Range("E5").Value = "'  -1.67%  "
$ws.Range("D6").Value = "'0.5008"
$ws.Range("E6").Value = "'  -3.00%  "
$ws.Range("D7").Value = "'1.010"
$ws.Range("E7").Value = "'  +0.21%  "
$ws.Range("D8").Value = "'0.2566"
$ws.Range("E8").Value = "'  +0.02%  "
$ws.Range("D9").Value = "'0.06412"
$ws.Range("E9").Value = "'  -0.43%  "
$ws.Range("D10").Value = "'19.44"
$ws.Range("E10").Value = "'  -2.67%  "
$ws.Range("D11").Value = "'0.07731"
$ws.Range("E11").Value = "'  +0.90%  "
$ws.Range("D12").Value = "'1.650.45"
$ws.Range("E12").Value = "'  -1.17%  "
$ws.Range("D13").Value = "'4.242"
$ws.Range("E13").Value = "'  -2.08%  "
$ws.Range("D14").Value = "'1.860.06"
$ws.Range("E14").Value = "'  -1.90%  "
$ws.Range("D15").Value = "'0.5430"
$ws.Range("E15").Value = "'  -2.03%  "
$ws.Range("D16").Value = "'0.0₅7901"
$ws.Range("E16").Value = "'  -1.72%  "
$ws.Range("D17").Value = "'63.40"
$ws.Range("E17").Value = "'  -1.86%  "
$ws.Range("D18").Value = "'25.865.20"
$ws.Range("E18").Value = "'  -1.94%  "
$ws.Range("D19").Value = "'1.010"
$ws.Range("E19").Value = "'  +0.34%  "
$ws.Range("D20").Value = "'203.29"
$ws.Range("E20").Value = "'  -3.40%  "
$ws.Range("D21").Value = "'4.289"
$ws.Range("E21").Value = "'  -2.71%  "
$ws.Range("D22").Value = "'9.971"
$ws.Range("E22").Value = "'  -1.38%  "
$ws.Range("D23").Value = "'5.925"
$ws.Range("E23").Value = "'  +0.58%  "
$ws.Range("D24").Value = "'1.011"
$ws.Range("E24").Value = "'  +0.23%  "
$ws.Range("D25").Value = "'1.958"
$ws.Range("E25").Value = "'  +12.25%  "
$ws.Range("D26").Value = "'140.77"
$ws.Range("E26").Value = "'  -2.75%  "
$ws.Range("D27").Value = "'0.1149"
$ws.Range("E27").Value = "'  -1.39%  "
$ws.Range("D28").Value = "'15.75"
$ws.Range("E28").Value = "'  -0.37%  "
$ws.Range("D29").Value = "'6.757"
$ws.Range("E29").Value = "'  -3.43%  "
$ws.Range("D30").Value = "'0.05080"
$ws.Range("E30").Value = "'  -3.41%  "
$ws.Range("D31").Value = "'1.239"
$ws.Range("E31").Value = "'  -1.94%  "
$ws.Range("D32").Value = "'3.252"
$ws.Range("E32").Value = "'  -3.57%  "
$ws.Range("D33").Value = "'3.185"
$ws.Range("E33").Value = "'  -1.03%  "
$ws.Range("D34").Value = "'1.540"
$ws.Range("E34").Value = "'  -2.14%  "
$ws.Range("D35").Value = "'2.343"
$ws.Range("E35").Value = "'  -1.45%  "
$ws.Range("D36").Value = "'0.8907"
$ws.Range("E36").Value = "'  -4.01%  "
$ws.Range("D37").Value = "'2.606"
$ws.Range("E37").Value = "'  -5.47%  "
$ws.Range("D38").Value = "'0.5632"
$ws.Range("E38").Value = "'  -1.65%  "
$ws.Range("D39").Value = "'1.136.31"
$ws.Range("E39").Value = "'  -0.89%  "
$ws.Range("D40").Value = "'0.01555"
$ws.Range("E40").Value = "'  -2.80%  "
$ws.Range("D41").Value = "'2.576"
$ws.Range("E41").Value = "'  -0.25%  "
$ws.Range("D42").Value = "'1.010"
$ws.Range("E42").Value = "'  +0.17%  "
$ws.Range("D43").Value = "'5.633"
$ws.Range("E43").Value = "'  -0.41%  "
$ws.Range("D44").Value = "'0.8158"
$ws.Range("E44").Value = "'  -3.57%  "
$ws.Range("D45").Value = "'99.30"
$ws.Range("E45").Value = "'  -0.73%  "
$ws.Range("D46").Value = "'1.770.18"
$ws.Range("E46").Value = "'  -2.03%  "
$ws.Range("E47").Value = "'  +1.58%  "
$ws.Range("D48").Value = "'0.4523"
$ws.Range("E48").Value = "'  +0.51%  "
$ws.Range("D50").Value = "'54.64"
$ws.Range("E50").Value = "'  -2.45%  "
$ws.Range("D51").Value = "'0.05020"
$ws.Range("E51").Value = "'  -1.67%  "

$ws.Range("D2:E51").ClearFormats()
